$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.210.80'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.579.79'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  -0.63%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.499'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.78%  '
$ws.Range('E7').Value = '  -0.57%  '
$ws.Range('E8').Value = '  -1.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.245'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.52'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.800.84'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.41%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.06'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.579.89'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.516'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.73%  '
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.194.77'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.99%  '
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.36'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '209.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.65%  '
$ws.Range('E24').Value = '  -2.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('E26').Value = '  -0.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('E28').Value = '  -2.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('E31').Value = '  -1.30%  '
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.283.92'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('E35').Value = '  -1.81%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.604'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.43%  '
$ws.Range('E37').Value = '  -1.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.12'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.74%  '
$ws.Range('E39').Value = '  -1.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.810'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.85%  '
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('E42').Value = '  +2.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.763'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.00%  '
$ws.Range('E44').Value = '  -3.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.34'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.713.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.48'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.93%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.99%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0102'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.101'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.03%  '
$ws.Range('E51').Value = '  -1.52%  '
